$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything currently in A:W
# (tip_documento .. estado) shifts one column to the right (B:X).
$ws.Range("A1").EntireColumn.Insert()

# Give the new header cell the same look as its neighbours (the other
# header cells all carry the bold/"header" cell style).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header text for the freshly inserted first column.
$ws.Range("A1").Value = "accion"

# Leave the selection on the new first cell.
$ws.Range("A1").Select() | Out-Null
